# Update the "cryptos" worksheet with refreshed Price (D) and Volume(1h) (E)
# values, matching the GitHub Actions data refresh.
# Values that look like plain numbers are prefixed with an apostrophe so
# Excel keeps them as text (preserving formatting such as trailing zeros
# and multi-dot "thousands" separators already used in this sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.517.52"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.594.23"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'207.27"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -4.21%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'22.20"
$ws.Range("E8").Value = "  -4.20%  "
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "'0.0587"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.822.93"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "1.591.08"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "'0.538"
$ws.Range("E15").Value = "  -3.59%  "
$ws.Range("D16").Value = "'63.31"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "27.510.23"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "'216.60"
$ws.Range("E18").Value = "  -5.06%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "0.0₃0689"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "'9.69"
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("D24").Value = "'2.00"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "'154.66"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "'15.00"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "1.348.15"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  -4.84%  "
$ws.Range("D43").Value = "'5.31"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").Value = "'63.77"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "'1.74"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").Value = "1.732.82"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("D48").Value = "'86.93"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "0.0₇0996"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").Value = "'0.0968"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "'0.0496"
$ws.Range("E51").Value = "  -1.19%  "
